# Commit: Tue, May 19, 2020 11:05:10 PM
#
# The underlying OOXML diff swaps the contents of ppt/theme/theme1.xml
# (the Slide Master's theme -- "Integral" / "Red Violet" color scheme)
# and ppt/theme/theme2.xml (the Notes Master's theme -- the stock
# "Office Theme" / "Office" color scheme). The font scheme and format
# scheme (fills/lines/effects) are byte-for-byte identical between the
# two themes in this deck, so the only structural difference between
# "Integral" and "Office Theme" is the 12-slot colour scheme. Re-create
# the post-edit theme1.xml colours (the "Office" palette) on the
# presentation's live theme via the SlideMaster's ThemeColorScheme.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Index order mirrors <a:clrScheme> child order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. RGB uses the usual VB encoding
# (r + g*256 + b*65536).
$cs.Item(1).RGB  = 0         # dk1      -> 000000
$cs.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      -> 44546A
$cs.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  -> FFC000
$cs.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$cs.Item(10).RGB = 4697456   # accent6  -> 70AD47
$cs.Item(11).RGB = 12673797  # hlink    -> 0563C1
$cs.Item(12).RGB = 7491477   # folHlink -> 954F72

Write-Output "Theme colour scheme updated to the Office Theme palette."
